# Update column F (dSF) values for specific rows, per repull/push of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    5  = 2
    6  = 2
    14 = -6
    15 = -4
    17 = 0
    18 = -2
    19 = -2
    27 = -5
    32 = 0
    36 = 2
    37 = 3
    44 = 1
    50 = 1
    52 = 15
    55 = -6
    59 = 6
    60 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
